$d = $word.ActiveDocument

# Locate the paragraph that currently reads
# "Pantalla de Vista Semanal/Mensual del Calendario" and the paragraph that
# currently reads the (older, duplicate) "Pantalla de Configuración" heading
# that follows the "Pantalla de Resumen de Productividad" section. Everything
# between them (both section bodies) gets removed, and the first heading's
# text is swapped in-place to become the new "Pantalla de Configuración"
# heading (keeping its run/paragraph formatting, including the
# lastRenderedPageBreak).

$startIndex = 0
$endIndex = 0

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($startIndex -eq 0 -and $t.StartsWith("Pantalla de Vista Semanal/Mensual del Calendario")) {
        $startIndex = $i
    }
    if ($startIndex -ne 0 -and $i -gt $startIndex -and $t.StartsWith("Pantalla de Configuraci")) {
        $endIndex = $i
        break
    }
}

# Delete the paragraphs strictly after the heading paragraph through the
# trailing duplicate "Pantalla de Configuración" heading (inclusive).
$delStart = $d.Paragraphs.Item($startIndex + 1)
$delEnd = $d.Paragraphs.Item($endIndex)
$delRange = $d.Range($delStart.Range.Start, $delEnd.Range.End)
$delRange.Delete()

# Rename the remaining heading paragraph's text in place. Using the
# Selection object (rather than Range.Text / Find.Execute) preserves the
# run's non-text children (e.g. <w:lastRenderedPageBreak/>) that a raw
# Range.Text reassignment would otherwise drop.
$headingRange = $d.Paragraphs.Item($startIndex).Range
$headingRange.Select()
$word.Selection.Text = "Pantalla de Configuración"
